$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (gp41) - row 10
$ws.Rows("10").Delete()

# Rename "Broad Spectrum" header to "EF50"
$ws.Range("C1").Value = "EF50"

# Add a new column D header "EF90" with the same style as the other headers
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "EF90"
$ws.Application.CutCopyMode = $false

# Update gRNA (B) and EF50 (C) counts, and fill in the new EF90 (D) counts
$ws.Range("B2").Value = 14
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 2

$ws.Range("B3").Value = 53
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = 135
$ws.Range("C4").Value = 82
$ws.Range("D4").Value = 2

$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0

$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 0

$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 0

$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0

$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
